$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- Refresh the re-query timestamps on the "data" sheet (col F, time_taken) ---
$data.Range("F2").Value = "2021-10-05 14:35:42.447667"
$data.Range("F3").Value = "2021-10-05 14:35:42.447675"
$data.Range("F4").Value = "2021-10-05 14:35:42.447678"
$data.Range("F5").Value = "2021-10-05 14:35:42.447680"
$data.Range("F6").Value = "2021-10-05 14:35:42.447683"
$data.Range("F7").Value = "2021-10-05 14:35:42.447686"
$data.Range("F8").Value = "2021-10-05 14:35:42.447689"
$data.Range("F9").Value = "2021-10-05 14:35:42.447691"
$data.Range("F10").Value = "2021-10-05 14:35:42.447694"

# --- Add the new "metadata" tab right after "data" ---
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Short Long Bones with Advanced Carpal Bone Age"
$meta.Range("C2").Value = 173
$meta.Range("D2").Value = "'0.1"
$meta.Range("E2").Value = "2020-01-15T06:45:27.256051Z"
$meta.Range("F2").Value = "2021-10-05 14:35:42.444134"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/173/?format=json"

# Match the bold / centered / thin-bordered header style used on the "data" sheet
$headerStyle = $meta.Range("B1:G1")
$headerStyle.Font.Bold = $true
$headerStyle.HorizontalAlignment = -4108
$headerStyle.VerticalAlignment = -4160
$headerStyle.Borders.LineStyle = 1

$meta.Range("A2").Font.Bold = $true
$meta.Range("A2").HorizontalAlignment = -4108
$meta.Range("A2").VerticalAlignment = -4160
$meta.Range("A2").Borders.LineStyle = 1

# Keep "data" as the active sheet (tab selection unchanged)
$data.Activate()
